# Apply the odds updates for the 2026-02-18 Betfair Back/Lay sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("L2").Value = 1.4
$ws.Range("G5").Value = 8.6
$ws.Range("H5").Value = 1.44
$ws.Range("I5").Value = 1.45
$ws.Range("J5").Value = 5.3
$ws.Range("N5").Value = 5.2
$ws.Range("R5").Value = 1.56
$ws.Range("T5").Value = 1.91
$ws.Range("U5").Value = 2.04
$ws.Range("V5").Value = 3.2
$ws.Range("W5").Value = 1.13
$ws.Range("Y5").Value = 9.6
$ws.Range("Z5").Value = 8.800000000000001
$ws.Range("AB5").Value = 32
$ws.Range("AJ5").Value = 270
$ws.Range("G6").Value = 5.2
$ws.Range("J6").Value = 4
$ws.Range("S6").Value = 2.36
$ws.Range("W6").Value = 1.24
$ws.Range("R7").Value = 1.49
$ws.Range("T7").Value = 1.72
$ws.Range("F9").Value = 2.34
$ws.Range("G9").Value = 2.36
$ws.Range("H9").Value = 3.5
$ws.Range("I9").Value = 3.55
$ws.Range("T9").Value = 1.81
$ws.Range("V9").Value = 1.39
$ws.Range("W9").Value = 1.73
$ws.Range("AA9").Value = 65
$ws.Range("AB9").Value = 10
$ws.Range("AI9").Value = 55
$ws.Range("F11").Value = 2.72
$ws.Range("G11").Value = 2.74
$ws.Range("H11").Value = 2.86
$ws.Range("I11").Value = 2.88
$ws.Range("W11").Value = 1.57
$ws.Range("Z11").Value = 19
$ws.Range("AK11").Value = 28
$ws.Range("H12").Value = 1.97
$ws.Range("I12").Value = 1.99
$ws.Range("O12").Value = 1.23
$ws.Range("V12").Value = 2
$ws.Range("AH12").Value = 15.5
$ws.Range("AJ12").Value = 80
$ws.Range("AO12").Value = 10.5
$ws.Range("H13").Value = 1.8
$ws.Range("I13").Value = 1.82
$ws.Range("J13").Value = 4.1
$ws.Range("K13").Value = 4.2
$ws.Range("R13").Value = 1.66
$ws.Range("S13").Value = 2.46
$ws.Range("V13").Value = 2.22
